# Backlog add artesanos api
#
# The "Productores de Artesanía Lenca" backlog rows for the new
# /api/artesanos endpoints (rows 4-8, column H = "Status") move from
# "Todo" to "Done" now that the artesanos API work is complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the artesanos API backlog items (rows 4-8) as Done.
$ws.Range("H4:H8").Value = "Done"

# Move the active selection to E9, matching where the editor left off.
$ws.Range("E9").Select()
